# Weekly cryptos list refresh (prices + 1h volume deltas).
# Source values are plain text (as scraped), so numeric-looking
# Price strings are forced to Text before assignment (NumberFormat="@")
# and the temporary number format is cleared right after, so the cell
# keeps its original (unstyled) look -- only .Value changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.847.33'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.09%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.461.40'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.18%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.72'
$ws.Range('D5').ClearFormats()

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.70%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.460.52'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.28%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.91%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.22'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.04%  '

# Row 11
$ws.Range('E11').Value = '  -1.70%  '

# Row 12
$ws.Range('E12').Value = '  -1.07%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.055.01'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.24%  '

# Row 14
$ws.Range('E14').Value = '  -0.01%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.60'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.07%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -9.94%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.874.82'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.14%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.472.02'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.91%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.77%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.45'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.15%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.99'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.58%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.04%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.51'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.01%  '

# Row 25
$ws.Range('E25').Value = '  -3.18%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000120'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.24%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.92'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.52%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.180'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.88%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.02%  '

# Row 30
$ws.Range('E30').Value = '  -3.22%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.48%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.02'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.22%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.32'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.27%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.06'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.23%  '

# Row 35
$ws.Range('E35').Value = '  -1.49%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '161.44'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.17%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.89'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.38%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.891.37'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.59%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0750'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.25%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.64'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.29%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.804'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.18%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.54'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.58%  '

# Row 43
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.56'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.26%  '

# Row 44
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.07'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.00%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.14'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.09%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0311'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.82%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.47'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +13.75%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '325.47'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.07%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.08'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.69%  '

# Row 50
$ws.Range('B50').Value = 'SuiNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.849'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.56%  '

# Row 51
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.47'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.17%  '

